$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "68.136.48"
$ws.Range("E2").Value = "  -0.15%  "

Set-TextValue "D3" "3.617.33"
$ws.Range("E3").Value = "  -1.45%  "

$ws.Range("E4").Value = "  -0.22%  "

Set-TextValue "D5" "587.50"
$ws.Range("E5").Value = "  -1.75%  "

Set-TextValue "D6" "193.68"
$ws.Range("E6").Value = "  +0.56%  "

Set-TextValue "D7" "3.611.09"
$ws.Range("E7").Value = "  -1.56%  "

$ws.Range("E8").Value = "  -0.59%  "

$ws.Range("E9").Value = "  +0.11%  "

Set-TextValue "D10" "0.681"
$ws.Range("E10").Value = "  -2.74%  "

$ws.Range("E11").Value = "  -1.07%  "

Set-TextValue "D12" "55.66"
$ws.Range("E12").Value = "  -2.87%  "

Set-TextValue "D13" "0.0000291"
$ws.Range("E13").Value = "  +6.59%  "

Set-TextValue "D14" "10.01"
$ws.Range("E14").Value = "  -2.45%  "

Set-TextValue "D15" "4.193.15"
$ws.Range("E15").Value = "  -1.60%  "

Set-TextValue "D16" "3.622.81"
$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("E17").Value = "  -0.45%  "

Set-TextValue "D18" "12.53"
$ws.Range("E18").Value = "  -0.63%  "

Set-TextValue "D19" "67.960.55"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("E20").Value = "  -2.06%  "

$ws.Range("E21").Value = "  -2.66%  "

Set-TextValue "D22" "405.51"
$ws.Range("E22").Value = "  -0.15%  "

Set-TextValue "D23" "13.58"
$ws.Range("E23").Value = "  +24.34%  "

Set-TextValue "D24" "4.28"
$ws.Range("E24").Value = "  -3.25%  "

Set-TextValue "D25" "86.12"
$ws.Range("E25").Value = "  -2.72%  "

Set-TextValue "D26" "2.96"
$ws.Range("E26").Value = "  +0.06%  "

Set-TextValue "D27" "12.66"
$ws.Range("E27").Value = "  +0.23%  "

Set-TextValue "D28" "3.91"
$ws.Range("E28").Value = "  +5.11%  "

$ws.Range("E29").Value = "  +0.71%  "

Set-TextValue "D30" "8.29"
$ws.Range("E30").Value = "  +15.37%  "

Set-TextValue "D31" "9.19"
$ws.Range("E31").Value = "  -2.02%  "

Set-TextValue "D32" "31.57"
$ws.Range("E32").Value = "  -1.39%  "

Set-TextValue "D33" "682.99"
$ws.Range("E33").Value = "  +11.90%  "

Set-TextValue "D34" "12.26"
$ws.Range("E34").Value = "  -0.53%  "

Set-TextValue "D35" "0.118"
$ws.Range("E35").Value = "  +1.31%  "

Set-TextValue "D36" "64.61"
$ws.Range("E36").Value = "  -3.82%  "

Set-TextValue "D37" "42.55"
$ws.Range("E37").Value = "  -3.75%  "

Set-TextValue "D38" "0.424"
$ws.Range("E38").Value = "  +7.56%  "

$ws.Range("E39").Value = "  +0.21%  "

Set-TextValue "D40" "0.0₃0789"
$ws.Range("E40").Value = "  +1.38%  "

Set-TextValue "D41" "3.00"
$ws.Range("E41").Value = "  +19.09%  "

Set-TextValue "D42" "3.15"
$ws.Range("E42").Value = "  +8.17%  "

Set-TextValue "D43" "3.199.48"
$ws.Range("E43").Value = "  +14.74%  "

$ws.Range("E44").Value = "  -1.32%  "

$ws.Range("E45").Value = "  -0.32%  "

Set-TextValue "D46" "0.0422"
$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("E47").Value = "  -2.86%  "

Set-TextValue "D48" "8.84"
$ws.Range("E48").Value = "  -1.22%  "

Set-TextValue "D49" "143.91"
$ws.Range("E49").Value = "  +0.31%  "

Set-TextValue "D50" "3.11"
$ws.Range("E50").Value = "  -3.94%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D51" "2.58"
$ws.Range("E51").Value = "  +1.75%  "
